$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EmptyString($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

# Common forecast values shared by every new row (7-15)
$commonB = 15.539999961853027
$commonC = 87.0
$commonD = 1014.0
$commonE = 1.5
$commonF = "2017-05-29T09:00:00"
$commonG = "2017-05-29T12:00:00"
$commonH = 14.079999923706055
$commonI = 995.1799926757812
$commonJ = 90.0
$commonPrecip = 0.3100000023841858

$timestamps = @{
    7  = "2017.05.29 03.41.21"
    8  = "2017.05.29 03.41.47"
    9  = "2017.05.29 03.42.21"
    10 = "2017.05.29 03.42.39"
    11 = "2017.05.29 03.43.42"
    12 = "2017.05.29 03.52.32"
    13 = "2017.05.29 03.58.20"
    14 = "2017.05.29 03.59.33"
    15 = "2017.05.29 03.59.45"
}

foreach ($row in 7..15) {
    $ws.Cells.Item($row, 1).Value = $timestamps[$row]
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $commonD
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    Set-EmptyString $row 11   # column K: "Forecasted Precipitation" left blank, as in source rows

    if ($row -eq 13) {
        # Row 13: precipitation value shifted out to column O, columns L-N blank
        Set-EmptyString $row 12
        Set-EmptyString $row 13
        Set-EmptyString $row 14
        $ws.Cells.Item($row, 15).Value = $commonPrecip
    }
    elseif ($row -eq 14 -or $row -eq 15) {
        # Rows 14-15: precipitation stays in column L, columns M-O blank
        $ws.Cells.Item($row, 12).Value = $commonPrecip
        Set-EmptyString $row 13
        Set-EmptyString $row 14
        Set-EmptyString $row 15
    }
    else {
        # Rows 7-12: precipitation in column L, no further columns
        $ws.Cells.Item($row, 12).Value = $commonPrecip
    }
}
